$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $escaped = $val.Replace('"', '""')
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy($ws.Range($addr)) | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null
}

$excel.CutCopyMode = 0

Set-TextValue $ws "D2" '245.28'
Set-TextValue $ws "D4" '5.222'
Set-TextValue $ws "D5" '0.05793'
Set-TextValue $ws "D8" '0.8152'
Set-TextValue $ws "D9" '0.8587'
Set-TextValue $ws "D10" '0.1362'
Set-TextValue $ws "D11" '0.06968'
Set-TextValue $ws "D12" '0.03191'
Set-TextValue $ws "D13" '0.02863'
Set-TextValue $ws "D14" '0.09372'
Set-TextValue $ws "D15" '3.742'
Set-TextValue $ws "D16" '0.001509'
Set-TextValue $ws "D17" '0.04705'
Set-TextValue $ws "B18" 'TigerCash'
Set-TextValue $ws "C18" 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws "D18" '0.006274'
Set-TextValue $ws "E18" '17TigerCashTCH'
Set-TextValue $ws "B19" 'BitKan'
Set-TextValue $ws "C19" 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue $ws "D19" '0.001237'
Set-TextValue $ws "E19" '18BitKanKAN'
Set-TextValue $ws "B20" 'HotbitToken'
Set-TextValue $ws "C20" 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue $ws "D20" '0.004537'
Set-TextValue $ws "E20" '19HotbitTokenHTB'
Set-TextValue $ws "B21" 'NitroEx'
Set-TextValue $ws "C21" 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-TextValue $ws "D21" '0.00006899'
Set-TextValue $ws "E21" '20NitroExNTX'
Set-TextValue $ws "B22" 'LEO'
Set-TextValue $ws "C22" 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws "D22" '3.502'
Set-TextValue $ws "E22" '21LEOLEO'
Set-TextValue $ws "B23" 'BTSEToken'
Set-TextValue $ws "C23" 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws "D23" '2.132'
Set-TextValue $ws "E23" '22BTSETokenBTSE'
Set-TextValue $ws "B24" 'One'
Set-TextValue $ws "C24" 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws "D24" '0.009788'
Set-TextValue $ws "E24" '23OneONE'
Set-TextValue $ws "D25" '0.3178'
Set-TextValue $ws "D26" '0.1355'
Set-TextValue $ws "D27" '0.1326'
Set-TextValue $ws "B41" 'KickToken'
Set-TextValue $ws "C41" 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws "D41" '0.006279'
Set-TextValue $ws "E41" '40KickTokenKICK'
Set-TextValue $ws "B42" 'BKEXToken'
Set-TextValue $ws "C42" 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws "D42" '0.1052'
Set-TextValue $ws "E42" '41BKEXTokenBKK'
Set-TextValue $ws "B43" 'CEJI'
Set-TextValue $ws "C43" 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue $ws "D43" '0.003400'
Set-TextValue $ws "E43" '42CEJICEJIBestin24h'
Set-TextValue $ws "D44" '0.007931'
Set-TextValue $ws "D45" '0.00005274'
Set-TextValue $ws "D48" '0.002345'

$excel.CutCopyMode = 0
